# Fruta / hortaliza, semanal
# Inserts a new weekly record at row 31, pushing the existing rows 31-64
# down to 32-65, and populates the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 31, shifting rows 31..64 down to 32..65
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly data point
$ws.Cells.Item(31, 1).Value = 1
$ws.Cells.Item(31, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(31, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(31, 4).Value = 44665
$ws.Cells.Item(31, 5).Value = 15
$ws.Cells.Item(31, 6).Value = 100112012
$ws.Cells.Item(31, 7).Value = "Espinaca"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 270
$ws.Cells.Item(31, 11).Value = 1800
$ws.Cells.Item(31, 12).Value = 2000
$ws.Cells.Item(31, 13).Value = 1900
$ws.Cells.Item(31, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(31, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 16).Value = 633
$ws.Cells.Item(31, 17).Value = 3
$ws.Cells.Item(31, 18).Value = "Hortaliza"
